# Adds rows 1992-2013 to Sheet1 (case data), fixing the "no case load" bug.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 1992
$ws.Cells.Item(1992, 1).Value = '21TRC05611'
$ws.Cells.Item(1992, 2).Value = 'Hemmeter'
$ws.Cells.Item(1992, 3).Value = 'OVI Alcohol / Drugs 3rd'
$ws.Cells.Item(1992, 4).Value = '4511.19A1A***'
$ws.Cells.Item(1992, 5).Value = 'UCM'
$ws.Cells.Item(1992, 6).Value = 'No Contest'
$ws.Cells.Item(1992, 7).Value = 'Guilty'
$ws.Cells.Item(1992, 8).NumberFormat = "@"
$ws.Cells.Item(1992, 8).Value = '$ 0'
$ws.Cells.Item(1992, 8).ClearFormats()
$ws.Cells.Item(1992, 9).NumberFormat = "@"
$ws.Cells.Item(1992, 9).Value = '$ 0'
$ws.Cells.Item(1992, 9).ClearFormats()
$ws.Cells.Item(1992, 10).NumberFormat = "@"
$ws.Cells.Item(1992, 10).Value = '5'
$ws.Cells.Item(1992, 10).ClearFormats()
$ws.Cells.Item(1992, 11).Value = 'None'

# Row 1993
$ws.Cells.Item(1993, 1).Value = '21TRC05611'
$ws.Cells.Item(1993, 2).Value = 'Hemmeter'
$ws.Cells.Item(1993, 3).Value = 'OVI Refusal 3rd/10yr Prior 20yr'
$ws.Cells.Item(1993, 4).Value = '4511.19A2***'
$ws.Cells.Item(1993, 5).Value = 'UCM'
$ws.Cells.Item(1993, 6).Value = 'No Contest'
$ws.Cells.Item(1993, 7).Value = 'Guilty'
$ws.Cells.Item(1993, 8).NumberFormat = "@"
$ws.Cells.Item(1993, 8).Value = '$ 0'
$ws.Cells.Item(1993, 8).ClearFormats()
$ws.Cells.Item(1993, 9).NumberFormat = "@"
$ws.Cells.Item(1993, 9).Value = '$ 0'
$ws.Cells.Item(1993, 9).ClearFormats()
$ws.Cells.Item(1993, 10).Value = 'None'
$ws.Cells.Item(1993, 11).Value = 'None'

# Row 1994
$ws.Cells.Item(1994, 1).Value = '21TRC05611'
$ws.Cells.Item(1994, 2).Value = 'Hemmeter'
$ws.Cells.Item(1994, 3).Value = 'Driving In Marked Lanes'
$ws.Cells.Item(1994, 4).NumberFormat = "@"
$ws.Cells.Item(1994, 4).Value = '4511.33'
$ws.Cells.Item(1994, 4).ClearFormats()
$ws.Cells.Item(1994, 5).Value = 'MM'
$ws.Cells.Item(1994, 6).Value = 'No Contest'
$ws.Cells.Item(1994, 7).Value = 'Guilty'
$ws.Cells.Item(1994, 8).NumberFormat = "@"
$ws.Cells.Item(1994, 8).Value = '$ 0'
$ws.Cells.Item(1994, 8).ClearFormats()
$ws.Cells.Item(1994, 9).NumberFormat = "@"
$ws.Cells.Item(1994, 9).Value = '$ 0'
$ws.Cells.Item(1994, 9).ClearFormats()
$ws.Cells.Item(1994, 10).Value = 'None'
$ws.Cells.Item(1994, 11).Value = 'None'

# Row 1995
$ws.Cells.Item(1995, 1).Value = '21TRC05611'
$ws.Cells.Item(1995, 2).Value = 'Hemmeter'
$ws.Cells.Item(1995, 3).Value = 'Seatbelt Required Driver'
$ws.Cells.Item(1995, 4).Value = '4513.263B1'
$ws.Cells.Item(1995, 5).Value = 'UCM'
$ws.Cells.Item(1995, 6).Value = 'No Contest'
$ws.Cells.Item(1995, 7).Value = 'Guilty'
$ws.Cells.Item(1995, 8).NumberFormat = "@"
$ws.Cells.Item(1995, 8).Value = '$ 0'
$ws.Cells.Item(1995, 8).ClearFormats()
$ws.Cells.Item(1995, 9).NumberFormat = "@"
$ws.Cells.Item(1995, 9).Value = '$ 0'
$ws.Cells.Item(1995, 9).ClearFormats()
$ws.Cells.Item(1995, 10).Value = 'None'
$ws.Cells.Item(1995, 11).Value = 'None'

# Row 1996
$ws.Cells.Item(1996, 1).Value = '21TRC05611'
$ws.Cells.Item(1996, 2).Value = 'Hemmeter'
$ws.Cells.Item(1996, 3).Value = 'OVI Alcohol / Drugs 3rd'
$ws.Cells.Item(1996, 4).Value = '4511.19A1A***'
$ws.Cells.Item(1996, 5).Value = 'UCM'
$ws.Cells.Item(1996, 6).Value = 'No Contest'
$ws.Cells.Item(1996, 7).Value = 'Guilty'
$ws.Cells.Item(1996, 8).NumberFormat = "@"
$ws.Cells.Item(1996, 8).Value = '$ 0'
$ws.Cells.Item(1996, 8).ClearFormats()
$ws.Cells.Item(1996, 9).NumberFormat = "@"
$ws.Cells.Item(1996, 9).Value = '$ 0'
$ws.Cells.Item(1996, 9).ClearFormats()
$ws.Cells.Item(1996, 10).Value = 'None'
$ws.Cells.Item(1996, 11).Value = 'None'

# Row 1997
$ws.Cells.Item(1997, 1).Value = '21TRC05611'
$ws.Cells.Item(1997, 2).Value = 'Hemmeter'
$ws.Cells.Item(1997, 3).Value = 'OVI Refusal 3rd/10yr Prior 20yr'
$ws.Cells.Item(1997, 4).Value = '4511.19A2***'
$ws.Cells.Item(1997, 5).Value = 'UCM'
$ws.Cells.Item(1997, 6).Value = 'No Contest'
$ws.Cells.Item(1997, 7).Value = 'Guilty'
$ws.Cells.Item(1997, 8).NumberFormat = "@"
$ws.Cells.Item(1997, 8).Value = '$ 0'
$ws.Cells.Item(1997, 8).ClearFormats()
$ws.Cells.Item(1997, 9).NumberFormat = "@"
$ws.Cells.Item(1997, 9).Value = '$ 0'
$ws.Cells.Item(1997, 9).ClearFormats()
$ws.Cells.Item(1997, 10).Value = 'None'
$ws.Cells.Item(1997, 11).Value = 'None'

# Row 1998
$ws.Cells.Item(1998, 1).Value = '21TRC05611'
$ws.Cells.Item(1998, 2).Value = 'Hemmeter'
$ws.Cells.Item(1998, 3).Value = 'Driving In Marked Lanes'
$ws.Cells.Item(1998, 4).NumberFormat = "@"
$ws.Cells.Item(1998, 4).Value = '4511.33'
$ws.Cells.Item(1998, 4).ClearFormats()
$ws.Cells.Item(1998, 5).Value = 'MM'
$ws.Cells.Item(1998, 6).Value = 'No Contest'
$ws.Cells.Item(1998, 7).Value = 'Guilty'
$ws.Cells.Item(1998, 8).NumberFormat = "@"
$ws.Cells.Item(1998, 8).Value = '$ 0'
$ws.Cells.Item(1998, 8).ClearFormats()
$ws.Cells.Item(1998, 9).NumberFormat = "@"
$ws.Cells.Item(1998, 9).Value = '$ 0'
$ws.Cells.Item(1998, 9).ClearFormats()
$ws.Cells.Item(1998, 10).Value = 'None'
$ws.Cells.Item(1998, 11).Value = 'None'

# Row 1999
$ws.Cells.Item(1999, 1).Value = '21TRC05611'
$ws.Cells.Item(1999, 2).Value = 'Hemmeter'
$ws.Cells.Item(1999, 3).Value = 'Seatbelt Required Driver'
$ws.Cells.Item(1999, 4).Value = '4513.263B1'
$ws.Cells.Item(1999, 5).Value = 'UCM'
$ws.Cells.Item(1999, 6).Value = 'No Contest'
$ws.Cells.Item(1999, 7).Value = 'Guilty'
$ws.Cells.Item(1999, 8).NumberFormat = "@"
$ws.Cells.Item(1999, 8).Value = '$ 0'
$ws.Cells.Item(1999, 8).ClearFormats()
$ws.Cells.Item(1999, 9).NumberFormat = "@"
$ws.Cells.Item(1999, 9).Value = '$ 0'
$ws.Cells.Item(1999, 9).ClearFormats()
$ws.Cells.Item(1999, 10).Value = 'None'
$ws.Cells.Item(1999, 11).Value = 'None'

# Row 2000
$ws.Cells.Item(2000, 1).Value = '21TRC05611'
$ws.Cells.Item(2000, 2).Value = 'Hemmeter'
$ws.Cells.Item(2000, 3).Value = 'OVI Alcohol / Drugs 3rd'
$ws.Cells.Item(2000, 4).Value = '4511.19A1A***'
$ws.Cells.Item(2000, 5).Value = 'UCM'
$ws.Cells.Item(2000, 6).Value = 'No Contest'
$ws.Cells.Item(2000, 7).Value = 'Guilty'
$ws.Cells.Item(2000, 8).NumberFormat = "@"
$ws.Cells.Item(2000, 8).Value = '$ 0'
$ws.Cells.Item(2000, 8).ClearFormats()
$ws.Cells.Item(2000, 9).NumberFormat = "@"
$ws.Cells.Item(2000, 9).Value = '$ 0'
$ws.Cells.Item(2000, 9).ClearFormats()
$ws.Cells.Item(2000, 10).Value = 'None'
$ws.Cells.Item(2000, 11).Value = 'None'

# Row 2001
$ws.Cells.Item(2001, 1).Value = '21TRC05611'
$ws.Cells.Item(2001, 2).Value = 'Hemmeter'
$ws.Cells.Item(2001, 3).Value = 'OVI Refusal 3rd/10yr Prior 20yr'
$ws.Cells.Item(2001, 4).Value = '4511.19A2***'
$ws.Cells.Item(2001, 5).Value = 'UCM'
$ws.Cells.Item(2001, 6).Value = 'No Contest'
$ws.Cells.Item(2001, 7).Value = 'Guilty'
$ws.Cells.Item(2001, 8).NumberFormat = "@"
$ws.Cells.Item(2001, 8).Value = '$ 0'
$ws.Cells.Item(2001, 8).ClearFormats()
$ws.Cells.Item(2001, 9).NumberFormat = "@"
$ws.Cells.Item(2001, 9).Value = '$ 0'
$ws.Cells.Item(2001, 9).ClearFormats()
$ws.Cells.Item(2001, 10).Value = 'None'
$ws.Cells.Item(2001, 11).Value = 'None'

# Row 2002
$ws.Cells.Item(2002, 1).Value = '21TRC05611'
$ws.Cells.Item(2002, 2).Value = 'Hemmeter'
$ws.Cells.Item(2002, 3).Value = 'Driving In Marked Lanes'
$ws.Cells.Item(2002, 4).NumberFormat = "@"
$ws.Cells.Item(2002, 4).Value = '4511.33'
$ws.Cells.Item(2002, 4).ClearFormats()
$ws.Cells.Item(2002, 5).Value = 'MM'
$ws.Cells.Item(2002, 6).Value = 'No Contest'
$ws.Cells.Item(2002, 7).Value = 'Guilty'
$ws.Cells.Item(2002, 8).NumberFormat = "@"
$ws.Cells.Item(2002, 8).Value = '$ 0'
$ws.Cells.Item(2002, 8).ClearFormats()
$ws.Cells.Item(2002, 9).NumberFormat = "@"
$ws.Cells.Item(2002, 9).Value = '$ 0'
$ws.Cells.Item(2002, 9).ClearFormats()
$ws.Cells.Item(2002, 10).Value = 'None'
$ws.Cells.Item(2002, 11).Value = 'None'

# Row 2003
$ws.Cells.Item(2003, 1).Value = '21TRC05611'
$ws.Cells.Item(2003, 2).Value = 'Hemmeter'
$ws.Cells.Item(2003, 3).Value = 'Seatbelt Required Driver'
$ws.Cells.Item(2003, 4).Value = '4513.263B1'
$ws.Cells.Item(2003, 5).Value = 'UCM'
$ws.Cells.Item(2003, 6).Value = 'No Contest'
$ws.Cells.Item(2003, 7).Value = 'Guilty'
$ws.Cells.Item(2003, 8).NumberFormat = "@"
$ws.Cells.Item(2003, 8).Value = '$ 0'
$ws.Cells.Item(2003, 8).ClearFormats()
$ws.Cells.Item(2003, 9).NumberFormat = "@"
$ws.Cells.Item(2003, 9).Value = '$ 0'
$ws.Cells.Item(2003, 9).ClearFormats()
$ws.Cells.Item(2003, 10).Value = 'None'
$ws.Cells.Item(2003, 11).Value = 'None'

# Row 2004
$ws.Cells.Item(2004, 1).Value = '21TRC10217'
$ws.Cells.Item(2004, 2).Value = 'Hemmeter'
$ws.Cells.Item(2004, 3).Value = 'OVI Alcohol / Drugs 1st'
$ws.Cells.Item(2004, 4).Value = '4511.19A1A*'
$ws.Cells.Item(2004, 5).Value = 'M1'
$ws.Cells.Item(2004, 6).Value = 'Guilty'
$ws.Cells.Item(2004, 7).Value = 'Guilty'
$ws.Cells.Item(2004, 8).NumberFormat = "@"
$ws.Cells.Item(2004, 8).Value = '$ 0'
$ws.Cells.Item(2004, 8).ClearFormats()
$ws.Cells.Item(2004, 9).NumberFormat = "@"
$ws.Cells.Item(2004, 9).Value = '$ 0'
$ws.Cells.Item(2004, 9).ClearFormats()
$ws.Cells.Item(2004, 10).Value = 'None'
$ws.Cells.Item(2004, 11).Value = 'None'

# Row 2005
$ws.Cells.Item(2005, 1).Value = '21TRC10217'
$ws.Cells.Item(2005, 2).Value = 'Hemmeter'
$ws.Cells.Item(2005, 3).Value = 'Turn And Stop Signals'
$ws.Cells.Item(2005, 4).NumberFormat = "@"
$ws.Cells.Item(2005, 4).Value = '4511.39'
$ws.Cells.Item(2005, 4).ClearFormats()
$ws.Cells.Item(2005, 5).Value = 'MM'
$ws.Cells.Item(2005, 6).Value = 'Guilty'
$ws.Cells.Item(2005, 7).Value = 'Guilty'
$ws.Cells.Item(2005, 8).NumberFormat = "@"
$ws.Cells.Item(2005, 8).Value = '$ 0'
$ws.Cells.Item(2005, 8).ClearFormats()
$ws.Cells.Item(2005, 9).NumberFormat = "@"
$ws.Cells.Item(2005, 9).Value = '$ 0'
$ws.Cells.Item(2005, 9).ClearFormats()
$ws.Cells.Item(2005, 10).Value = 'None'
$ws.Cells.Item(2005, 11).Value = 'None'

# Row 2006
$ws.Cells.Item(2006, 1).Value = '21TRC10217 22TRD1234'
$ws.Cells.Item(2006, 2).Value = 'Hemmeter'
$ws.Cells.Item(2006, 3).Value = 'OVI Alcohol / Drugs 1st'
$ws.Cells.Item(2006, 4).Value = '4511.19A1A*'
$ws.Cells.Item(2006, 5).Value = 'M1'
$ws.Cells.Item(2006, 6).Value = 'Not Guilty'

# Row 2007
$ws.Cells.Item(2007, 1).Value = '21TRC10217 22TRD1234'
$ws.Cells.Item(2007, 2).Value = 'Hemmeter'
$ws.Cells.Item(2007, 3).Value = 'Turn And Stop Signals'
$ws.Cells.Item(2007, 4).NumberFormat = "@"
$ws.Cells.Item(2007, 4).Value = '4511.39'
$ws.Cells.Item(2007, 4).ClearFormats()
$ws.Cells.Item(2007, 5).Value = 'MM'
$ws.Cells.Item(2007, 6).Value = 'Not Guilty'

# Row 2008
$ws.Cells.Item(2008, 1).Value = '21TRC10217 22TRD1234'
$ws.Cells.Item(2008, 2).Value = 'Hemmeter'
$ws.Cells.Item(2008, 3).Value = 'Criminal Mischief M3'
$ws.Cells.Item(2008, 4).NumberFormat = "@"
$ws.Cells.Item(2008, 4).Value = '2909.07'
$ws.Cells.Item(2008, 4).ClearFormats()
$ws.Cells.Item(2008, 5).Value = 'M3'
$ws.Cells.Item(2008, 6).Value = 'Not Guilty'

# Row 2009
$ws.Cells.Item(2009, 1).Value = '21TRC05611'
$ws.Cells.Item(2009, 2).Value = 'Hemmeter'
$ws.Cells.Item(2009, 3).Value = 'OVI Alcohol / Drugs 3rd'
$ws.Cells.Item(2009, 4).Value = '4511.19A1A***'
$ws.Cells.Item(2009, 5).Value = 'UCM'
$ws.Cells.Item(2009, 6).Value = 'No Contest'
$ws.Cells.Item(2009, 7).Value = 'Guilty'
$ws.Cells.Item(2009, 8).NumberFormat = "@"
$ws.Cells.Item(2009, 8).Value = '$ 0'
$ws.Cells.Item(2009, 8).ClearFormats()
$ws.Cells.Item(2009, 9).NumberFormat = "@"
$ws.Cells.Item(2009, 9).Value = '$ 0'
$ws.Cells.Item(2009, 9).ClearFormats()
$ws.Cells.Item(2009, 10).Value = 'None'
$ws.Cells.Item(2009, 11).Value = 'None'

# Row 2010
$ws.Cells.Item(2010, 1).Value = '21TRC05611'
$ws.Cells.Item(2010, 2).Value = 'Hemmeter'
$ws.Cells.Item(2010, 3).Value = 'OVI Refusal 3rd/10yr Prior 20yr'
$ws.Cells.Item(2010, 4).Value = '4511.19A2***'
$ws.Cells.Item(2010, 5).Value = 'UCM'
$ws.Cells.Item(2010, 6).Value = 'No Contest'
$ws.Cells.Item(2010, 7).Value = 'Guilty'
$ws.Cells.Item(2010, 8).NumberFormat = "@"
$ws.Cells.Item(2010, 8).Value = '$ 0'
$ws.Cells.Item(2010, 8).ClearFormats()
$ws.Cells.Item(2010, 9).NumberFormat = "@"
$ws.Cells.Item(2010, 9).Value = '$ 0'
$ws.Cells.Item(2010, 9).ClearFormats()
$ws.Cells.Item(2010, 10).Value = 'None'
$ws.Cells.Item(2010, 11).Value = 'None'

# Row 2011
$ws.Cells.Item(2011, 1).Value = '21TRC05611'
$ws.Cells.Item(2011, 2).Value = 'Hemmeter'
$ws.Cells.Item(2011, 3).Value = 'Driving In Marked Lanes'
$ws.Cells.Item(2011, 4).NumberFormat = "@"
$ws.Cells.Item(2011, 4).Value = '4511.33'
$ws.Cells.Item(2011, 4).ClearFormats()
$ws.Cells.Item(2011, 5).Value = 'MM'
$ws.Cells.Item(2011, 6).Value = 'No Contest'
$ws.Cells.Item(2011, 7).Value = 'Guilty'
$ws.Cells.Item(2011, 8).NumberFormat = "@"
$ws.Cells.Item(2011, 8).Value = '$ 0'
$ws.Cells.Item(2011, 8).ClearFormats()
$ws.Cells.Item(2011, 9).NumberFormat = "@"
$ws.Cells.Item(2011, 9).Value = '$ 0'
$ws.Cells.Item(2011, 9).ClearFormats()
$ws.Cells.Item(2011, 10).Value = 'None'
$ws.Cells.Item(2011, 11).Value = 'None'

# Row 2012
$ws.Cells.Item(2012, 1).Value = '21TRC05611'
$ws.Cells.Item(2012, 2).Value = 'Hemmeter'
$ws.Cells.Item(2012, 3).Value = 'Seatbelt Required Driver'
$ws.Cells.Item(2012, 4).Value = '4513.263B1'
$ws.Cells.Item(2012, 5).Value = 'UCM'
$ws.Cells.Item(2012, 6).Value = 'No Contest'
$ws.Cells.Item(2012, 7).Value = 'Guilty'
$ws.Cells.Item(2012, 8).NumberFormat = "@"
$ws.Cells.Item(2012, 8).Value = '$ 0'
$ws.Cells.Item(2012, 8).ClearFormats()
$ws.Cells.Item(2012, 9).NumberFormat = "@"
$ws.Cells.Item(2012, 9).Value = '$ 0'
$ws.Cells.Item(2012, 9).ClearFormats()
$ws.Cells.Item(2012, 10).Value = 'None'
$ws.Cells.Item(2012, 11).Value = 'None'

# Row 2013
$ws.Cells.Item(2013, 1).Value = 'ASDF'
$ws.Cells.Item(2013, 2).Value = 'Bunner'
$ws.Cells.Item(2013, 3).Value = 'Disorderly Conduct - Persistent'
$ws.Cells.Item(2013, 4).Value = '2917.11(A)(1)'
$ws.Cells.Item(2013, 5).Value = 'M4'
$ws.Cells.Item(2013, 6).Value = 'No Contest'
$ws.Cells.Item(2013, 7).Value = 'Guilty'
$ws.Cells.Item(2013, 8).NumberFormat = "@"
$ws.Cells.Item(2013, 8).Value = '$ 0'
$ws.Cells.Item(2013, 8).ClearFormats()
$ws.Cells.Item(2013, 9).NumberFormat = "@"
$ws.Cells.Item(2013, 9).Value = '$ 0'
$ws.Cells.Item(2013, 9).ClearFormats()
